$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1166.0834
$ws.Range("I18").Value = 1181.1818
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1181.1818
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -897.1818000000001
$ws.Range("N18").Value = -1568
$ws.Range("H106").Value = 2332.2307
$ws.Range("I106").Value = 1732
$ws.Range("J106").Value = 4333
$ws.Range("K106").Value = 1732
$ws.Range("L106").Value = 4333
$ws.Range("M106").Value = -1101
$ws.Range("N106").Value = -5595
$ws.Range("H125").Value = 840702.1
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 840702.1
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 7566318.899999999
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -7571238.899999999
$ws.Range("H130").Value = 12141.429
$ws.Range("J130").Value = 12141.429
$ws.Range("L130").Value = 12141.429
$ws.Range("N130").Value = -22181.429
$ws.Range("H131").Value = 4441.73
$ws.Range("I131").Value = 736.38464
$ws.Range("J131").Value = 4995.4023
$ws.Range("K131").Value = 2209.15392
$ws.Range("L131").Value = 14986.2069
$ws.Range("M131").Value = 2830.84608
$ws.Range("N131").Value = -25066.2069
$ws.Range("H135").Value = 2694.4
$ws.Range("I135").Value = 1200
$ws.Range("J135").Value = 3690.6667
$ws.Range("K135").Value = 10800
$ws.Range("L135").Value = 33216.0003
$ws.Range("M135").Value = -8265
$ws.Range("N135").Value = -38286.0003
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 1343.8108
$ws.Range("I137").Value = 1424.7778
$ws.Range("J137").Value = 1267.1052
$ws.Range("K137").Value = 4274.3334
$ws.Range("L137").Value = 3801.3156
$ws.Range("M137").Value = -1724.3334
$ws.Range("N137").Value = -8901.3156
$ws.Range("H138").Value = 14928217
$ws.Range("I138").Value = 1204.86
$ws.Range("J138").Value = 58831190
$ws.Range("K138").Value = 3614.58
$ws.Range("L138").Value = 176493570
$ws.Range("M138").Value = 1525.42
$ws.Range("N138").Value = -176503850
$ws.Range("H141").Value = 3634.3215
$ws.Range("I141").Value = 1298.5385
$ws.Range("J141").Value = 33999.5
$ws.Range("K141").Value = 3895.6155
$ws.Range("L141").Value = 101998.5
$ws.Range("M141").Value = 1284.3845
$ws.Range("N141").Value = -112358.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 49514.5
$ws.Range("I33").Value = 29000
$ws.Range("K33").Value = 29000
$ws.Range("M33").Value = -28671
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H74").Value = 1209.6595
$ws.Range("I74").Value = 1164.0488
$ws.Range("K74").Value = 1164.0488
$ws.Range("M74").Value = -290.0488
$ws.Range("H77").Value = 1209.6595
$ws.Range("I77").Value = 1164.0488
$ws.Range("K77").Value = 5820.244000000001
$ws.Range("M77").Value = -1452.244000000001
$ws.Range("H130").Value = 38990
$ws.Range("J130").Value = 38990
$ws.Range("L130").Value = 38990
$ws.Range("N130").Value = -49030
$ws.Range("H131").Value = 63154.6
$ws.Range("J131").Value = 63154.6
$ws.Range("L131").Value = 63154.6
$ws.Range("N131").Value = -73234.60000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H135").Value = 44680
$ws.Range("J135").Value = 44680
$ws.Range("L135").Value = 44680
$ws.Range("N135").Value = -54820
$ws.Range("H140").Value = 60949
$ws.Range("J140").Value = 60949
$ws.Range("L140").Value = 60949
$ws.Range("N140").Value = -71309

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 288432.38
$ws.Range("J5").Value = 439695.72
$ws.Range("L5").Value = 1319087.16
$ws.Range("N5").Value = -1319311.16
$ws.Range("H113").Value = 1596.2
$ws.Range("I113").Value = 2027.75
$ws.Range("J113").Value = 1308.5
$ws.Range("K113").Value = 6083.25
$ws.Range("L113").Value = 3925.5
$ws.Range("M113").Value = -3913.25
$ws.Range("N113").Value = -8265.5
$ws.Range("H135").Value = 288432.38
$ws.Range("J135").Value = 439695.72
$ws.Range("L135").Value = 3957261.48
$ws.Range("N135").Value = -3962331.48
$ws.Range("H137").Value = 2714.1
$ws.Range("I137").Value = 1802.579
$ws.Range("K137").Value = 5407.737
$ws.Range("M137").Value = -307.7370000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 509.5
$ws.Range("I17").Value = 119
$ws.Range("J17").Value = 900
$ws.Range("K17").Value = 119
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = 49
$ws.Range("N17").Value = -1236
$ws.Range("H41").Value = 1518.3334
$ws.Range("I41").Value = 777.5
$ws.Range("K41").Value = 777.5
$ws.Range("M41").Value = -422.5
$ws.Range("H140").Value = 49996.668
$ws.Range("J140").Value = 49996.668
$ws.Range("L140").Value = 49996.668
$ws.Range("N140").Value = -60356.668

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62378.766
$ws.Range("I7").Value = 86953.336
$ws.Range("J7").Value = 3399.8
$ws.Range("K7").Value = 86953.336
$ws.Range("L7").Value = 3399.8
$ws.Range("M7").Value = -86841.336
$ws.Range("N7").Value = -3623.8
$ws.Range("H22").Value = 1065.3
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 1575.75
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 1575.75
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -2165.75
$ws.Range("H25").Value = 18966.666
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 27950
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 27950
$ws.Range("M25").Value = -770
$ws.Range("N25").Value = -28410
$ws.Range("H27").Value = 1065.3
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 1575.75
$ws.Range("K27").Value = 725
$ws.Range("L27").Value = 1575.75
$ws.Range("M27").Value = -618
$ws.Range("N27").Value = -1789.75
$ws.Range("H40").Value = 25212.436
$ws.Range("I40").Value = 31899.705
$ws.Range("K40").Value = 31899.705
$ws.Range("M40").Value = -31763.705
$ws.Range("H93").Value = 100003
$ws.Range("I93").Value = 100003
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 100003
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -98755
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 62378.766
$ws.Range("I126").Value = 86953.336
$ws.Range("J126").Value = 3399.8
$ws.Range("K126").Value = 260860.008
$ws.Range("L126").Value = 10199.4
$ws.Range("M126").Value = -258390.008
$ws.Range("N126").Value = -15139.4
$ws.Range("H130").Value = 39500
$ws.Range("J130").Value = 39500
$ws.Range("L130").Value = 39500
$ws.Range("N130").Value = -49540
$ws.Range("H136").Value = 22227050
$ws.Range("I136").Value = 5458.2856
$ws.Range("J136").Value = 41670944
$ws.Range("K136").Value = 16374.8568
$ws.Range("L136").Value = 125012832
$ws.Range("M136").Value = -13824.8568
$ws.Range("N136").Value = -125017932

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 63292
$ws.Range("J130").Value = 63292
$ws.Range("L130").Value = 63292
$ws.Range("N130").Value = -73332
$ws.Range("H135").Value = 72739.625
$ws.Range("J135").Value = 72739.625
$ws.Range("L135").Value = 72739.625
$ws.Range("N135").Value = -82879.625
$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 50000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 50000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -60360
